$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-14 06:52:21"

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
